# Update the NMOP & BBF "Meeting Objectives" slide (slide 4):
#  - refine the title line
#  - rewrite the "Collaboration..." sentence with emphasis runs
#  - rename "Note:" to "LSes Log:"
#  - update/annotate the first liaison-statement record (color, size, hyperlink)
#  - append a second liaison-statement record (new paragraph + hyperlink)
#  - append a trailing empty paragraph
#  - nudge the placeholder's position

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# ---------------------------------------------------------------------------
# 1. Paragraph 1 - title sentence
# ---------------------------------------------------------------------------
$full = $tr.Text
$old1 = "NMOP & BBF exchange of information "
$idx = $full.IndexOf($old1)
$sub = $tr.Characters($idx + 1, $old1.Length)
$sub.Text = "NMOP & BBF Initial exchange of information"

# ---------------------------------------------------------------------------
# 2. Paragraph 2 - "Collaboration ... industry?" -> long sentence w/ emphasis
# ---------------------------------------------------------------------------
$full = $tr.Text
$old2 = "Collaboration a consistent data collection architecture for the industry?"
$idx = $full.IndexOf($old2)
$sub = $tr.Characters($idx + 1, $old2.Length)
$sub.Text = "Collaboration to ensure a consistent data collection architecture for the industry by leveraging common building blocks"

# Re-locate and style the sub-segments that need bold+italic+blue.
$full = $tr.Text

function Style-Accent($range) {
    $range.Font.Bold = $true
    $range.Font.Italic = $true
    $range.Font.Color.RGB = 15773696  # 0x00B0F0 (RRGGBB) -> BGR int
}

$needle = "ensure"
$idx = $full.IndexOf($needle)
Style-Accent ($tr.Characters($idx + 1, $needle.Length))

$needle = " a consistent data collection architecture "
$idx = $full.IndexOf($needle)
Style-Accent ($tr.Characters($idx + 1, $needle.Length))

$needle = "leveraging "
$idx = $full.IndexOf($needle)
Style-Accent ($tr.Characters($idx + 1, $needle.Length))

$needle = "common"
$idx = $full.IndexOf($needle)
Style-Accent ($tr.Characters($idx + 1, $needle.Length))

$needle = " building blocks"
$idx = $full.IndexOf($needle)
Style-Accent ($tr.Characters($idx + 1, $needle.Length))

# ---------------------------------------------------------------------------
# 3. Paragraph 4 - "Note:" -> "LSes Log:"
# ---------------------------------------------------------------------------
$full = $tr.Text
$old4 = "Note:"
$idx = $full.IndexOf($old4)
$sub = $tr.Characters($idx + 1, $old4.Length)
$sub.Text = "LSes Log:"

# ---------------------------------------------------------------------------
# 4. Paragraph 5 - first LS record: punctuation + color/size + hyperlink
# ---------------------------------------------------------------------------
$full = $tr.Text
$old5a = "IETF NMOP LS to BBF on Automated Intelligent Management (AIM), WT-508: Broadband Network Data Collection (BNDC)"
$idx = $full.IndexOf($old5a)
$sub = $tr.Characters($idx + 1, $old5a.Length)
$quote1 = [char]0x201C
$quote2 = [char]0x201D
$sub.Text = "IETF NMOP LS to BBF on " + $quote1 + "Automated Intelligent Management (AIM), WT-508: Broadband Network Data Collection (BNDC)" + $quote2 + ": "
$sub.Font.Size = 12
$sub.Font.Color.RGB = 6299648  # 0x002060 (RRGGBB) -> BGR int

# The line break run between the title text and the URL becomes plain text
# (the hard break disappears, replaced by the new trailing runs).
$full = $tr.Text
$oldUrl1 = "https://datatracker.ietf.org/liaison/1969/"
$idx = $full.IndexOf($oldUrl1)
$sub = $tr.Characters($idx + 1, $oldUrl1.Length)
$sub.Font.Size = 12
$sub.ActionSettings.Item(1).Hyperlink.Address = $oldUrl1
$sub.InsertAfter(" (09/12/2024)")

# ---------------------------------------------------------------------------
# 5. New paragraph (6) - second LS record, appended after paragraph 5
# ---------------------------------------------------------------------------
$full = $tr.Text
$idx = $full.IndexOf($oldUrl1 + " (09/12/2024)")
$endOfPara5 = $idx + ($oldUrl1 + " (09/12/2024)").Length
$anchor = $tr.Characters($endOfPara5, 1)
$newPara2Text = "Response to IETF NMOP on " + $quote1 + "Automated Intelligent Management (AIM)and Broadband Network Data Collection (BNDC)" + $quote2 + ": " + "https://datatracker.ietf.org/liaison/1975/" + "  (31/01/25)"
$anchor.InsertAfter("`r" + $newPara2Text)

# Style the new paragraph's runs (all sz=1200, bold; URL gets hyperlink).
$full = $tr.Text
$urlneedle = "https://datatracker.ietf.org/liaison/1975/"
$idxUrl = $full.IndexOf($urlneedle)
$preText = "Response to IETF NMOP on " + $quote1 + "Automated Intelligent Management (AIM)and Broadband Network Data Collection (BNDC)" + $quote2 + ": "
$idxPre = $full.IndexOf($preText)
$subPre = $tr.Characters($idxPre + 1, $preText.Length)
$subPre.Font.Size = 12
$subPre.Font.Bold = $true

$subUrl = $tr.Characters($idxUrl + 1, $urlneedle.Length)
$subUrl.Font.Size = 12
$subUrl.Font.Bold = $true
$subUrl.ActionSettings.Item(1).Hyperlink.Address = $urlneedle

$tailNeedle = "  (31/01/25)"
$idxTail = $full.LastIndexOf($tailNeedle)
$subTail = $tr.Characters($idxTail + 1, $tailNeedle.Length)
$subTail.Font.Size = 12
$subTail.Font.Bold = $true

# ---------------------------------------------------------------------------
# 6. Trailing empty paragraph (level 2 / no bullet)
# ---------------------------------------------------------------------------
$full = $tr.Text
$idxTailEnd = $full.LastIndexOf($tailNeedle) + $tailNeedle.Length
$anchor2 = $tr.Characters($idxTailEnd, 1)
$anchor2.InsertAfter("`r")
$lastParas = $tr.Paragraphs()
$lastPara = $lastParas.Item($lastParas.Count)
$lastPara.Font.Size = 12
$lastPara.Font.Bold = $true
$lastPara.IndentLevel = 2
$lastPara.ParagraphFormat.Bullet.Visible = $false

# ---------------------------------------------------------------------------
# 7. Shape position nudge
# ---------------------------------------------------------------------------
$sh.Left = 627489 / 12700.0
$sh.Top = 1723929 / 12700.0
